# Regenerate orders with updated distance/sizes.
#
# The workbook encodes a trial-order table (Trial, Condition, Duration_Seconds,
# Filename_Left, Filename_Right, Is_Repeat, Block, Distance, Face, Size,
# ConditionID). Several stimulus parameters were renamed:
#   Distance: D64 -> D69, D51 -> D55, D80 -> D86
#   Size:     S30 -> S31
# These tokens show up both standalone (Distance / Size columns) and embedded
# inside composite labels (Condition / Filename_Left / Filename_Right, e.g.
# "Face18_D64_S20" or "Face18_D64_S20_l.png"), so every text cell in the used
# range gets the same substring substitutions applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @(
    @("D64", "D69"),
    @("D51", "D55"),
    @("D80", "D86"),
    @("S30", "S31")
)

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$rowOffset = $used.Row
$colOffset = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = $val
            foreach ($pair in $replacements) {
                $newVal = $newVal.Replace($pair[0], $pair[1])
            }
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
